$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply updated values for rows 2-11 across columns A-T (natmi recomputed stats
# plus a new "M2" target-cluster category inserted into the permutation).

# Row 2
$ws.Range("A2").Value = "FAPs"
$ws.Range("B2").Value = "Angpt4"
$ws.Range("C2").Value = "Tek"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 1.027368666666667
$ws.Range("H2").Value = 3.082106
$ws.Range("I2").Value = 0.7595351715900243
$ws.Range("J2").Value = 0.8257207729145039
$ws.Range("K2").Value = 2
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 70.2375145
$ws.Range("N2").Value = 140.475029
$ws.Range("O2").Value = 0.5875234044920492
$ws.Range("P2").Value = 0.4906443493732269
$ws.Range("Q2").Value = 72.15982162184567
$ws.Range("R2").Value = 432.958929731074
$ws.Range("S2").Value = 0.4462446898440238
$ws.Range("T2").Value = 0.4051352313905948

# Row 3
$ws.Range("A3").Value = "FAPs"
$ws.Range("B3").Value = "Angpt4"
$ws.Range("C3").Value = "Tek"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 1.027368666666667
$ws.Range("H3").Value = 3.082106
$ws.Range("I3").Value = 0.7595351715900243
$ws.Range("J3").Value = 0.8257207729145039
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 2.804800666666666
$ws.Range("N3").Value = 8.414401999999999
$ws.Range("O3").Value = 0.02346162230159168
$ws.Range("P3").Value = 0.02938941407625393
$ws.Range("Q3").Value = 2.881564321179111
$ws.Range("R3").Value = 25.93407889061199
$ws.Range("S3").Value = 0.01781992732061978
$ws.Range("T3").Value = 0.0242674497065488

# Row 4
$ws.Range("A4").Value = "FAPs"
$ws.Range("B4").Value = "Angpt4"
$ws.Range("C4").Value = "Tek"
$ws.Range("D4").Value = "Neutro"
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 1.027368666666667
$ws.Range("H4").Value = 3.082106
$ws.Range("I4").Value = 0.7595351715900243
$ws.Range("J4").Value = 0.8257207729145039
$ws.Range("K4").Value = 1
$ws.Range("L4").Value = 0.3333333333333333
$ws.Range("M4").Value = 0.1062743333333333
$ws.Range("N4").Value = 0.318823
$ws.Range("O4").Value = 0.0008889645166775211
$ws.Range("P4").Value = 0.001113569468636453
$ws.Range("Q4").Value = 0.1091829201375556
$ws.Range("R4").Value = 0.9826462812379999
$ws.Range("S4").Value = 0.000675199816712104
$ws.Range("T4").Value = 0.0009194974423364854

# Row 5
$ws.Range("A5").Value = "FAPs"
$ws.Range("B5").Value = "Angpt4"
$ws.Range("C5").Value = "Tek"
$ws.Range("D5").Value = "M2"
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 1.027368666666667
$ws.Range("H5").Value = 3.082106
$ws.Range("I5").Value = 0.7595351715900243
$ws.Range("J5").Value = 0.8257207729145039
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 44.29925666666667
$ws.Range("N5").Value = 132.89777
$ws.Range("O5").Value = 0.3705548278372964
$ws.Range("P5").Value = 0.4641788676534301
$ws.Range("Q5").Value = 45.51166825595778
$ws.Range("R5").Value = 409.60501430362
$ws.Range("S5").Value = 0.2814494247449129
$ws.Range("T5").Value = 0.3832821333693695

# Row 6
$ws.Range("A6").Value = "FAPs"
$ws.Range("B6").Value = "Angpt4"
$ws.Range("C6").Value = "Tek"
$ws.Range("D6").Value = "sCs"
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 1.027368666666667
$ws.Range("H6").Value = 3.082106
$ws.Range("I6").Value = 0.7595351715900243
$ws.Range("J6").Value = 0.8257207729145039
$ws.Range("K6").Value = 2
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 2.1006075
$ws.Range("N6").Value = 4.201215
$ws.Range("O6").Value = 0.01757118085238527
$ws.Range("P6").Value = 0.01467379942845245
$ws.Range("Q6").Value = 2.158098326465
$ws.Range("R6").Value = 12.94858995879
$ws.Range("S6").Value = 0.01334592986375579
$ws.Range("T6").Value = 0.01211646100565416

# Row 7
$ws.Range("A7").Value = "sCs"
$ws.Range("B7").Value = "Angpt4"
$ws.Range("C7").Value = "Tek"
$ws.Range("D7").Value = "ECs"
$ws.Range("E7").Value = 2
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 0.3252595
$ws.Range("H7").Value = 0.6505190000000001
$ws.Range("I7").Value = 0.2404648284099757
$ws.Range("J7").Value = 0.1742792270854962
$ws.Range("K7").Value = 2
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 70.2375145
$ws.Range("N7").Value = 140.475029
$ws.Range("O7").Value = 0.5875234044920492
$ws.Range("P7").Value = 0.4906443493732269
$ws.Range("Q7").Value = 22.84541884751275
$ws.Range("R7").Value = 91.38167539005102
$ws.Range("S7").Value = 0.1412787146480254
$ws.Range("T7").Value = 0.08550911798263214

# Row 8
$ws.Range("A8").Value = "sCs"
$ws.Range("B8").Value = "Angpt4"
$ws.Range("C8").Value = "Tek"
$ws.Range("D8").Value = "FAPs"
$ws.Range("E8").Value = 2
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 0.3252595
$ws.Range("H8").Value = 0.6505190000000001
$ws.Range("I8").Value = 0.2404648284099757
$ws.Range("J8").Value = 0.1742792270854962
$ws.Range("K8").Value = 3
$ws.Range("L8").Value = 1
$ws.Range("M8").Value = 2.804800666666666
$ws.Range("N8").Value = 8.414401999999999
$ws.Range("O8").Value = 0.02346162230159168
$ws.Range("P8").Value = 0.02938941407625393
$ws.Range("Q8").Value = 0.9122880624396668
$ws.Range("R8").Value = 5.473728374638
$ws.Range("S8").Value = 0.005641694980971904
$ws.Range("T8").Value = 0.005121964369705136

# Row 9
$ws.Range("A9").Value = "sCs"
$ws.Range("B9").Value = "Angpt4"
$ws.Range("C9").Value = "Tek"
$ws.Range("D9").Value = "Neutro"
$ws.Range("E9").Value = 2
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 0.3252595
$ws.Range("H9").Value = 0.6505190000000001
$ws.Range("I9").Value = 0.2404648284099757
$ws.Range("J9").Value = 0.1742792270854962
$ws.Range("K9").Value = 1
$ws.Range("L9").Value = 0.3333333333333333
$ws.Range("M9").Value = 0.1062743333333333
$ws.Range("N9").Value = 0.318823
$ws.Range("O9").Value = 0.0008889645166775211
$ws.Range("P9").Value = 0.001113569468636453
$ws.Range("Q9").Value = 0.03456673652283334
$ws.Range("R9").Value = 0.207400419137
$ws.Range("S9").Value = 0.0002137646999654171
$ws.Range("T9").Value = 0.0001940720262999677

# Row 10
$ws.Range("A10").Value = "sCs"
$ws.Range("B10").Value = "Angpt4"
$ws.Range("C10").Value = "Tek"
$ws.Range("D10").Value = "M2"
$ws.Range("E10").Value = 2
$ws.Range("F10").Value = 1
$ws.Range("G10").Value = 0.3252595
$ws.Range("H10").Value = 0.6505190000000001
$ws.Range("I10").Value = 0.2404648284099757
$ws.Range("J10").Value = 0.1742792270854962
$ws.Range("K10").Value = 3
$ws.Range("L10").Value = 1
$ws.Range("M10").Value = 44.29925666666667
$ws.Range("N10").Value = 132.89777
$ws.Range("O10").Value = 0.3705548278372964
$ws.Range("P10").Value = 0.4641788676534301
$ws.Range("Q10").Value = 14.40875407377167
$ws.Range("R10").Value = 86.45252444263001
$ws.Range("S10").Value = 0.08910540309238357
$ws.Range("T10").Value = 0.08089673428406062

# Row 11
$ws.Range("A11").Value = "sCs"
$ws.Range("B11").Value = "Angpt4"
$ws.Range("C11").Value = "Tek"
$ws.Range("D11").Value = "sCs"
$ws.Range("E11").Value = 2
$ws.Range("F11").Value = 1
$ws.Range("G11").Value = 0.3252595
$ws.Range("H11").Value = 0.6505190000000001
$ws.Range("I11").Value = 0.2404648284099757
$ws.Range("J11").Value = 0.1742792270854962
$ws.Range("K11").Value = 2
$ws.Range("L11").Value = 1
$ws.Range("M11").Value = 2.1006075
$ws.Range("N11").Value = 4.201215
$ws.Range("O11").Value = 0.01757118085238527
$ws.Range("P11").Value = 0.01467379942845245
$ws.Range("Q11").Value = 0.6832425451462502
$ws.Range("R11").Value = 2.732970180585001
$ws.Range("S11").Value = 0.004225250988629473
$ws.Range("T11").Value = 0.002557338422798289

